$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Update column C ("Förändrad") for every data row (2-36): 46065 -> 46066
for ($row = 2; $row -le 36; $row++) {
    $ws.Cells.Item($row, 3).Value = 46066
}

# Step 2: Apply the row permutation for rows 10-36 (columns A, B, G).
# This re-orders the underlying records; row 22 keeps its original content.
# Capture original values first so the row-by-row writes do not clobber sources still needed later.
$origA = @{}
$origB = @{}
$origG = @{}
for ($row = 10; $row -le 36; $row++) {
    $origA[$row] = $ws.Cells.Item($row, 1).Value2
    $origB[$row] = $ws.Cells.Item($row, 2).Value2
    $origG[$row] = $ws.Cells.Item($row, 7).Value2
}

$rowMap = @{
    10 = 11
    11 = 10
    12 = 23
    13 = 27
    14 = 30
    15 = 35
    16 = 31
    17 = 20
    18 = 25
    19 = 18
    20 = 14
    21 = 12
    22 = 22
    23 = 26
    24 = 13
    25 = 29
    26 = 33
    27 = 34
    28 = 32
    29 = 24
    30 = 36
    31 = 15
    32 = 16
    33 = 17
    34 = 21
    35 = 19
    36 = 28
}

foreach ($row in $rowMap.Keys) {
    $src = $rowMap[$row]
    $ws.Cells.Item($row, 1).Value = $origA[$src]
    $ws.Cells.Item($row, 2).Value = $origB[$src]
    $ws.Cells.Item($row, 7).Value = $origG[$src]
}
